# Insert a new price record as row 265 in the data table, shifting all
# subsequent rows (old 265-371) down by one (new 266-372).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(265).Insert()

$ws.Range("A265").Value = 3
$ws.Range("B265").Value = "Femacal de La Calera"
$ws.Range("C265").Value = "Coquimbo"
$ws.Range("D265").Value = 45176
$ws.Range("E265").Value = 5
$ws.Range("F265").Value = "Fruta"
$ws.Range("G265").Value = 100101
$ws.Range("H265").Value = "Berries"
$ws.Range("I265").Value = 100101001
$ws.Range("J265").Value = "Arándano (blue)"
$ws.Range("K265").Value = "Sin especificar"
$ws.Range("L265").Value = "Primera"
$ws.Range("M265").Value = 40
$ws.Range("N265").Value = 13000
$ws.Range("O265").Value = 13000
$ws.Range("P265").Value = 13000
$ws.Range("Q265").Value = "$/bandeja 12 canastillos 125 gramos"
$ws.Range("R265").Value = "Provincia de Quillota"
$ws.Range("S265").Value = 8667
$ws.Range("T265").Value = 1.5
